$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the SQL in B5: remove the redundant CONCAT(...) wrapper around REPLACE(...)
$old = "CONCAT(REPLACE(trt.treatment_agent, ';', ', ')) AS ""Treatment Agent"","
$new = "REPLACE(trt.treatment_agent, ';', ', ') AS ""Treatment Agent"","

$cell = $ws.Range("B5")
$text = $cell.Value2
$text = $text.Replace($old, $new)
$cell.Value = $text

# Update the sheet view: drop the frozen/scrolled topLeftCell and move the
# active selection from C5 to B2 (also resets the visible top-left cell to A1).
$ws.Range("B2").Select()
